# Update the "想去人数" (interest count) values in column F across the
# workbook's sheets, as produced by the latest gh-pages data generation run.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 8088
$ws1.Range("F4").Value  = 93
$ws1.Range("F5").Value  = 30688
$ws1.Range("F15").Value = 387
$ws1.Range("F19").Value = 413
$ws1.Range("F24").Value = 2335
$ws1.Range("F25").Value = 828
$ws1.Range("F30").Value = 1074

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 343

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value  = 522

# Sheet "全部类型" (All Types, aggregate of the above)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 522
$ws4.Range("F3").Value  = 8088
$ws4.Range("F5").Value  = 93
$ws4.Range("F7").Value  = 30688
$ws4.Range("F15").Value = 343
$ws4.Range("F21").Value = 387
$ws4.Range("F29").Value = 413
$ws4.Range("F34").Value = 2335
$ws4.Range("F35").Value = 828
$ws4.Range("F41").Value = 1074
